$wb = $excel.ActiveWorkbook
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
Write-Host $wb.Worksheets.Count
